# Updates cryptos list values (price + volume%) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.187.44'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.834.98'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''241.39'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").Value = '''0.6656'
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.07413'
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").Value = '''0.2935'
$ws.Range("E9").Value = '  -2.29%  '
$ws.Range("D10").Value = '''22.59'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '''0.07721'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '1.857.90'
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("D13").Value = '''4.992'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '''0.6686'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = '''83.02'
$ws.Range("E15").Value = '  -5.66%  '
$ws.Range("D16").Value = '''6.101'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '''0.000008410'
$ws.Range("E17").Value = '  +2.54%  '
$ws.Range("D18").Value = '29.222.59'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").Value = '''226.48'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").Value = '''12.46'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").Value = '''7.180'
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '''159.61'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = '''0.1414'
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''8.623'
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("D27").Value = '''17.94'
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("D28").Value = '''1.513'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").Value = '''4.109'
$ws.Range("E29").Value = '  -4.08%  '
$ws.Range("D30").Value = '''4.046'
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("D31").Value = '''1.183'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").Value = '''0.05308'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7611'
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '''1.875'
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '1.272.06'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").Value = '''0.01796'
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("D39").Value = '''2.722'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = '''0.9296'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("D41").Value = '''0.09099'
$ws.Range("E41").Value = '  +19.64%  '
$ws.Range("D42").Value = '''5.959'
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").Value = '''1.003'
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").Value = '''102.77'
$ws.Range("D45").Value = '2.017.02'
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = '''0.5163'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").Value = '''1.777'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("D49").Value = '''63.36'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").Value = '''0.05923'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").Value = '''8.920'
$ws.Range("E51").Value = '  -5.71%  '
